$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "survey" sheet
#   - add two new trailing columns (E: appearance, F: parameters)
#   - rewrite rows 2-3 with new question definitions
#
# Columns B/C originally hold a genuinely-blank cell (an empty text value,
# not "no cell") in rows 2-3 (the option list's blank placeholders). A plain
# value assignment of "" always clears a cell outright in this host, so to
# keep that "blank text" cell alive in its new column E we relocate the
# whole column structurally (Insert/Delete), which moves existing cells
# instead of re-writing them, and only overwrite real content afterwards.
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# Shift columns B:D three slots to the right (-> E:G), carrying column B's
# blank cells (rows 2-3) along with them so they land in column E.
$survey.Columns("B:D").Insert()
# Drop the old "label"/"required" columns that rode along in F:G - they'll
# be re-created fresh below.
$survey.Columns("F:G").Delete()

$survey.Range("A1").Value = "type"
$survey.Range("B1").Value = "name"
$survey.Range("C1").Value = "label"
$survey.Range("D1").Value = "required"
$survey.Range("E1").Value = "appearance"
$survey.Range("F1").Value = "parameters"

$survey.Range("A2").Value = "text"
$survey.Range("B2").Value = "what_is_your_name?"
$survey.Range("C2").Value = "What is your name?"
$survey.Range("D2").Value = $false

$survey.Range("A3").Value = "integer"
$survey.Range("B3").Value = "what_is_your_age?"
$survey.Range("C3").Value = "What is your age?"
$survey.Range("D3").Value = $false

# ---------------------------------------------------------------------------
# "choices" sheet: the zEdBNsf option list is gone, drop rows 2-3 entirely.
# ---------------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")
$choices.Rows("2:3").Delete()
